# Natmi following Dr Hou advice
# Update the Vwf-Itga2b LR-pairs sheet with recomputed statistics
# (ligand/receptor-expressing cell counts changed from 1 to 3, and all
# dependent expression / specificity / edge-weight values recomputed).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2 = @{ 'E' = 3; 'G' = 28.67276066666667; 'H' = 86.018282; 'I' = 0.9474462168692853; 'J' = 0.9474462168692853; 'K' = 3; 'M' = 1.104136666666667; 'N' = 3.31241; 'O' = 0.2772362398998524; 'P' = 0.2772362398998525; 'Q' = 31.65864638662444; 'R' = 284.92781747962; 'S' = 0.2626664266721808; 'T' = 0.2626664266721808 }
    3 = @{ 'E' = 3; 'G' = 28.67276066666667; 'H' = 86.018282; 'I' = 0.9474462168692853; 'J' = 0.9474462168692853; 'K' = 3; 'M' = 2.230986666666666; 'N' = 6.692959999999999; 'O' = 0.5601755411317187; 'P' = 0.5601755411317187; 'Q' = 63.96854674385777; 'R' = 575.71692069472; 'S' = 0.5307361972279516; 'T' = 0.5307361972279516 }
    4 = @{ 'E' = 3; 'G' = 28.67276066666667; 'H' = 86.018282; 'I' = 0.9474462168692853; 'J' = 0.9474462168692853; 'K' = 3; 'M' = 0.647533; 'N' = 1.942599; 'O' = 0.1625882189684289; 'P' = 0.1625882189684289; 'Q' = 18.56655873276867; 'R' = 167.099028594918; 'S' = 0.1540435929691529; 'T' = 0.1540435929691529 }
    5 = @{ 'E' = 3; 'G' = 0.7569533333333333; 'H' = 2.27086; 'I' = 0.02501233070476559; 'J' = 0.02501233070476559; 'K' = 3; 'M' = 1.104136666666667; 'N' = 3.31241; 'O' = 0.2772362398998524; 'P' = 0.2772362398998525; 'Q' = 0.8357799302888888; 'R' = 7.522019372599999; 'S' = 0.006934324515720837; 'T' = 0.006934324515720839 }
    6 = @{ 'E' = 3; 'G' = 0.7569533333333333; 'H' = 2.27086; 'I' = 0.02501233070476559; 'J' = 0.02501233070476559; 'K' = 3; 'M' = 2.230986666666666; 'N' = 6.692959999999999; 'O' = 0.5601755411317187; 'P' = 0.5601755411317187; 'Q' = 1.688752793955555; 'R' = 15.1987751456; 'S' = 0.01401129588750757; 'T' = 0.01401129588750757 }
    7 = @{ 'E' = 3; 'G' = 0.7569533333333333; 'H' = 2.27086; 'I' = 0.02501233070476559; 'J' = 0.02501233070476559; 'K' = 3; 'M' = 0.647533; 'N' = 1.942599; 'O' = 0.1625882189684289; 'P' = 0.1625882189684289; 'Q' = 0.4901522627933333; 'R' = 4.41137036514; 'S' = 0.004066710301537184; 'T' = 0.004066710301537184 }
    8 = @{ 'E' = 3; 'G' = 0.8334926666666668; 'H' = 2.500478; 'I' = 0.02754145242594914; 'J' = 0.02754145242594913; 'K' = 3; 'M' = 1.104136666666667; 'N' = 3.31241; 'O' = 0.2772362398998524; 'P' = 0.2772362398998525; 'Q' = 0.9202898146644445; 'R' = 8.28260833198; 'S' = 0.007635488711950808; 'T' = 0.007635488711950808 }
    9 = @{ 'E' = 3; 'G' = 0.8334926666666668; 'H' = 2.500478; 'I' = 0.02754145242594914; 'J' = 0.02754145242594913; 'K' = 3; 'M' = 2.230986666666666; 'N' = 6.692959999999999; 'O' = 0.5601755411317187; 'P' = 0.5601755411317187; 'Q' = 1.859511026097778; 'R' = 16.73559923488; 'S' = 0.01542804801625955; 'T' = 0.01542804801625954 }
    10 = @{ 'E' = 3; 'G' = 0.8334926666666668; 'H' = 2.500478; 'I' = 0.02754145242594914; 'J' = 0.02754145242594913; 'K' = 3; 'M' = 0.647533; 'N' = 1.942599; 'O' = 0.1625882189684289; 'P' = 0.1625882189684289; 'Q' = 0.5397140069246668; 'R' = 4.857426062322; 'S' = 0.004477915697738785; 'T' = 0.004477915697738784 }
}

foreach ($rowNum in $data.Keys) {
    $rowData = $data[$rowNum]
    foreach ($col in $rowData.Keys) {
        $ws.Range("$col$rowNum").Value = $rowData[$col]
    }
}
